$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "Gaizka"
$ws.Range("C12").Value = "Tareas en proceso del listado "
$ws.Range("D12").Value = $ws.Range("D11").Value2 + 2

$ws.Range("D12").NumberFormat = $ws.Range("D11").NumberFormat

$ws.Range("D13").Select()
